# Update odds/stats values on Sheet1 per the 2024-11-05 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 3 (San Lorenzo - Estudiantes L.P.)
$ws.Range("G3").Value = 2.9
$ws.Range("I3").Value = 2.8
$ws.Range("AA3").Value = 41
$ws.Range("AI3").Value = 11

# Row 7 (Botafogo RJ - Vasco)
$ws.Range("O7").Value = 1.22
$ws.Range("P7").Value = 4.33
$ws.Range("Q7").Value = 1.73
$ws.Range("R7").Value = 2.08

# Row 8 (Internacional - Criciuma)
$ws.Range("M8").Value = 1.03
$ws.Range("N8").Value = 15
$ws.Range("Q8").Value = 1.67
$ws.Range("R8").Value = 2.15

# Row 9 (Amazonas - America MG)
$ws.Range("G9").Value = 3.1
$ws.Range("I9").Value = 2.38
$ws.Range("J9").Value = 3.75
$ws.Range("K9").Value = 2.05
$ws.Range("L9").Value = 3.1
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8
$ws.Range("R9").Value = 1.67
$ws.Range("S9").Value = 1.44
$ws.Range("T9").Value = 2.63
$ws.Range("W9").Value = 9
$ws.Range("X9").Value = 15
$ws.Range("Y9").Value = 12
$ws.Range("Z9").Value = 34
$ws.Range("AC9").Value = 8
$ws.Range("AI9").Value = 11
$ws.Range("AK9").Value = 23
$ws.Range("AL9").Value = 21
$ws.Range("AN9").Value = 5
$ws.Range("AT9").Value = 2.63
$ws.Range("AW9").Value = 4.33
$ws.Range("AX9").Value = 13
$ws.Range("AZ9").Value = 41
$ws.Range("BA9").Value = 67

# Row 10 (Brusque - Botafogo SP)
$ws.Range("Q10").Value = 2.88
$ws.Range("R10").Value = 1.4

# Row 13 (Club Leon - Mazatlan FC)
$ws.Range("Q13").Value = 1.7
$ws.Range("R13").Value = 2.1

# Row 14 (Santos Laguna - Guadalajara Chivas)
$ws.Range("G14").Value = 5.75
$ws.Range("I14").Value = 1.53
$ws.Range("N14").Value = 12
$ws.Range("Z14").Value = 67
$ws.Range("AC14").Value = 12
$ws.Range("AN14").Value = 7.5
$ws.Range("AR14").Value = 126
$ws.Range("AU14").Value = 8.5
$ws.Range("AY14").Value = 19
